$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" (H) column values to complete the consolidated report.
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
